# Applies updated cryptocurrency Price (column D) and Volume(1h) (column E)
# values scraped on Sat Feb  4 04:57:47 UTC 2023, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format ("@") before writing so values such as "329.50" or
# "2.09%" are stored as literal text (matching the source data) instead of
# being auto-converted by Excel into numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "329.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.09%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.82%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.633"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-4.27%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08161"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.75%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.041"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "6.23%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.752"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.10%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.528"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.99%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.949"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.03%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9182"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.49%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1255"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.56%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1954"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.14%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09414"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.33%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03714"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "4.90%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1055"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "10.34%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001308"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.86%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006149"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.93%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.437"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.52%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-2.11%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.258"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-5.30%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.80%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2651"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "10.05%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04423"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.31%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001268"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.55%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.67%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.67%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "13.39%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05429"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.86%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007650"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.63%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009434"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.61%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.87%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.31%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01285"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "15.83%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006873"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.93%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.12%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002283"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "60.48%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003524"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "17.40%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.12%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.12%"
